$p = $ppt.ActivePresentation

# --- Update cached "datetimeFigureOut" field text across master / layouts / notes master ---
# The fields were regenerated from 12/5/2019 (US) / 5-12-2019 (NL) to 12/15/2019 / 15-12-2019.

# Slide master date placeholder (nl-NL, "5-12-2019" -> "15-12-2019")
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "5-12-2019") {
            $sh.TextFrame.TextRange.Text = "15-12-2019"
        }
    }
}

# All slide layouts under the master share the same date placeholder text
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "5-12-2019") {
                $sh.TextFrame.TextRange.Text = "15-12-2019"
            }
        }
    }
}

# Notes master date placeholder (en-US, "12/5/2019" -> "12/15/2019")
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "12/5/2019") {
            $sh.TextFrame.TextRange.Text = "12/15/2019"
        }
    }
}

# --- Exercise timing text trimmed off the two "calculator" slides ---

# Slide 22 title: "Calculator 25 min " -> "Calculator"
$s22 = $p.Slides.Item(22)
for ($i = 1; $i -le $s22.Shapes.Count; $i++) {
    $sh = $s22.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "Calculator 25 min ") {
            $sh.TextFrame.TextRange.Text = "Calculator"
        }
    }
}

# Slide 24 title: "Assignment build the calculator (25 min)" -> "Assignment build the calculator"
$s24 = $p.Slides.Item(24)
for ($i = 1; $i -le $s24.Shapes.Count; $i++) {
    $sh = $s24.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "Assignment build the calculator (25 min)") {
            $sh.TextFrame.TextRange.Text = "Assignment build the calculator"
        }
    }
}
